# Update handback/handoff timestamps for the zh-cn and de-de report rows.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 corresponds to the 583a8f0c... entry.
$wsZhCn.Range("E2").Value = "2016-03-14 01:21:09"
$wsZhCn.Range("H2").Value = "2016-03-14 01:21:39"

# de-de sheet: row 2 corresponds to the 583a8f0c... entry.
$wsDeDe.Range("E2").Value = "2016-03-14 01:21:13"
$wsDeDe.Range("H2").Value = "2016-03-14 01:21:46"
